$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts existing rows 2-77 down to 3-78)
$ws.Rows("2:2").Insert()

# Populate the new row with the "id" translation key entry
$ws.Range("A2").Value = "id"
$ws.Range("B2").Value = "id"
$ws.Range("C2").Value = "?"

# The inserted row picked up the bold header style from row 1 above it;
# reset the font weight back to normal to match the rest of the data rows.
$ws.Range("A2:C2").Font.Bold = $false

# Move the selection/cursor to B7 and scroll the view back to the top
$ws.Range("B7").Select() | Out-Null
